$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "v61"
$ws.Range("A7").Value = "v71"
$ws.Range("A9").Value = "v91"
$ws.Range("A11").Value = "v11_1"
$ws.Range("B12").Value = "v12_2"

$ws.Range("B13").Select() | Out-Null
